$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: update Version, Date and Contact values
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# ---------------------------------------------------------------------------
# 2) "Include from FSIII" sheet: insert two new concept rows (FBOE) right
#    after the header row, pushing the existing C / D / (blank) / System URI
#    rows down, then fill in the new concept codes and restore formatting.
# ---------------------------------------------------------------------------
$inc = $wb.Worksheets.Item("Include from FSIII")

# Shift the existing concept rows (rows 2-5) down by two rows.
$inc.Range("A2:B3").Insert(-4121)  # xlShiftDown

# Re-apply the standard data-row formatting (copied from the row right below,
# which still carries the original style) to the two freshly inserted rows.
$inc.Range("A4:B4").Copy()
$inc.Range("A2:B3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the two new concept codes.
$inc.Range("A2").Value = "d6d48a71-b96f-4b88-86f9-b13bd3c03560"
$inc.Range("A3").Value = "687159ad-a61c-47c0-a878-53aa54bae2d5"
